# HDX Design Challenge - add subtitle line to the title, and tweak the
# Piloting-the-Pilot caption on slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title 1 (shape 1): grow the title box upward/taller and add a
#     second line ("Creating a self-modifying digital workflow in 2 days")
#     below the existing "HDX Design Challenge" title, at a smaller size.
$title = $s.Shapes.Item(1)

# Resize/reposition (only Top/Height change per the target; Left/Width are
# untouched so we leave them alone to avoid introducing any drift).
$title.Top = 0.4239370078740157
$title.Height = 111.5659842519685

$titleTr = $title.TextFrame.TextRange
$titleTr.Text = "HDX Design Challenge" + [char]13 + "Creating a self-modifying digital workflow in 2 days"

# Size the second line (the new subtitle) at 31pt.
$subtitlePara = $titleTr.Paragraphs(2, 1)
$subtitlePara.Font.Size = 31

# --- TextBox 24 (shape 21): widen the box slightly and reword the caption.
$caption = $s.Shapes.Item(21)
$caption.Width = 386.2656792913386
$caption.TextFrame.TextRange.Text = "HDX 2021 (Piloting the Pilot) Certification Program"
